$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row (row 1): shift the mortgage/rent headers right by one column
# and add two new trailing headers.
$ws.Range("K1").Value = "Mortgage Interest Rate"
$ws.Range("L1").Value = "Mortgage Payoff Duration"
$ws.Range("M1").Value = "Monthly Rent"
$ws.Range("N1").Value = "HomeLocation"
$ws.Range("O1").Value = "DownPayment"

# --- Row 3 ("Framework_002" / Mortgage Calc test case) ---
# A3 test-case name changes from Framework_002 to Mortgage Calc
$ws.Range("A3").Value = "Mortgage Calc"

# J3 used to hold a phone number string; now holds the numeric home price
$ws.Range("J3").Value = 80000

# K3 used to be a hyperlinked e-mail address; drop the hyperlink and store
# the numeric interest rate instead
$ws.Range("K3").Hyperlinks.Delete()
$ws.Range("K3").Value = 12

# New trailing columns for this row
$ws.Range("N3").Value = "Ben Claire, SD"
$ws.Range("O3").Value = 20000

# --- Row 4 ("Buy vs Rent" test case) ---
# The monthly-rent value moves from column N to column M, and N4 is cleared
$ws.Range("M4").Value = 987879
$ws.Range("N4").ClearContents()
